$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 107
$ws.Range("J9").Value = 36.666668
$ws.Range("L9").Value = 36.666668
$ws.Range("N9").Value = -374.666668

$ws.Range("H28").Value = 811
$ws.Range("I28").Value = 588.9167
$ws.Range("J28").Value = 1699.3334
$ws.Range("K28").Value = 588.9167
$ws.Range("L28").Value = 1699.3334
$ws.Range("M28").Value = -103.9167
$ws.Range("N28").Value = -2669.3334

$ws.Range("H70").Value = 38323.6
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 246990.67
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 740972.01
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -741512.01

$ws.Range("H73").Value = 38323.6
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 246990.67
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 740972.01
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -742844.01

$ws.Range("H74").Value = 3501.1428
$ws.Range("I74").Value = 2901.8
$ws.Range("J74").Value = 4999.5
$ws.Range("K74").Value = 2901.8
$ws.Range("L74").Value = 4999.5
$ws.Range("M74").Value = -1965.8
$ws.Range("N74").Value = -6871.5

$ws.Range("H77").Value = 3501.1428
$ws.Range("I77").Value = 2901.8
$ws.Range("J77").Value = 4999.5
$ws.Range("K77").Value = 14509
$ws.Range("L77").Value = 24997.5
$ws.Range("M77").Value = -9829
$ws.Range("N77").Value = -34357.5

$ws.Range("H129").Value = 1809
$ws.Range("I129").Value = 557.1429000000001
$ws.Range("K129").Value = 1671.4287
$ws.Range("M129").Value = 3328.5713

$ws.Range("H141").Value = 1521.4333
$ws.Range("I141").Value = 1521.4333
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4564.2999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 615.7001
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 30000
$ws.Range("J10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("N10").Value = -30340

$ws.Range("H138").Value = 74979.664
$ws.Range("J138").Value = 74979.664
$ws.Range("L138").Value = 74979.664
$ws.Range("N138").Value = -85259.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 7500
$ws.Range("I96").Value = 7500
$ws.Range("K96").Value = 7500
$ws.Range("M96").Value = -4754

$ws.Range("H105").Value = 3526.6365
$ws.Range("I105").Value = 3310.4443
$ws.Range("J105").Value = 4499.5
$ws.Range("K105").Value = 3310.4443
$ws.Range("L105").Value = 4499.5
$ws.Range("M105").Value = -1563.4443
$ws.Range("N105").Value = -7993.5

$ws.Range("H107").Value = 8666.666999999999
$ws.Range("I107").Value = 8666.666999999999
$ws.Range("K107").Value = 8666.666999999999
$ws.Range("M107").Value = -6746.666999999999

$ws.Range("H134").Value = 3733.2778
$ws.Range("I134").Value = 3199.8333
$ws.Range("K134").Value = 9599.499899999999
$ws.Range("M134").Value = -7064.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 700
$ws.Range("I5").Value = 700
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 700
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -588
$ws.Range("N5").ClearContents()

$ws.Range("H7").Value = 152.07692
$ws.Range("I7").Value = 61.625
$ws.Range("K7").Value = 61.625
$ws.Range("M7").Value = 51.375

$ws.Range("H22").Value = 405.46155
$ws.Range("I22").Value = 247.5
$ws.Range("J22").Value = 434.18182
$ws.Range("K22").Value = 247.5
$ws.Range("L22").Value = 434.18182
$ws.Range("M22").Value = 102.5
$ws.Range("N22").Value = -1134.18182

$ws.Range("H105").Value = 1304.0555
$ws.Range("I105").Value = 980.63635
$ws.Range("J105").Value = 1812.2858
$ws.Range("K105").Value = 980.63635
$ws.Range("L105").Value = 1812.2858
$ws.Range("M105").Value = 766.36365
$ws.Range("N105").Value = -5306.2858

$ws.Range("H122").Value = 2640.05
$ws.Range("I122").Value = 2695.4443
$ws.Range("K122").Value = 8086.3329
$ws.Range("M122").Value = -5636.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1624.7
$ws.Range("I5").Value = 814.8570999999999
$ws.Range("J5").Value = 3514.3333
$ws.Range("K5").Value = 2444.5713
$ws.Range("L5").Value = 10542.9999
$ws.Range("M5").Value = -2332.5713
$ws.Range("N5").Value = -10766.9999

$ws.Range("H12").Value = 263.16666
$ws.Range("I12").Value = 223.22223
$ws.Range("K12").Value = 669.66669
$ws.Range("M12").Value = -496.66669

$ws.Range("H41").Value = 200
$ws.Range("I41").Value = 200
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 600
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -262
$ws.Range("N41").ClearContents()

$ws.Range("H68").Value = 1516.7333
$ws.Range("I68").Value = 1100
$ws.Range("J68").Value = 2141.8333
$ws.Range("K68").Value = 3300
$ws.Range("L68").Value = 6425.499899999999
$ws.Range("M68").Value = -2489
$ws.Range("N68").Value = -8047.499899999999

$ws.Range("H71").Value = 1516.7333
$ws.Range("I71").Value = 1100
$ws.Range("J71").Value = 2141.8333
$ws.Range("K71").Value = 9900
$ws.Range("L71").Value = 19276.4997
$ws.Range("M71").Value = -5844
$ws.Range("N71").Value = -27388.4997

$ws.Range("H135").Value = 1624.7
$ws.Range("I135").Value = 814.8570999999999
$ws.Range("J135").Value = 3514.3333
$ws.Range("K135").Value = 7333.7139
$ws.Range("L135").Value = 31628.9997
$ws.Range("M135").Value = -4798.7139
$ws.Range("N135").Value = -36698.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H122").Value = 41054.77
$ws.Range("I122").Value = 2164.5715
$ws.Range("J122").Value = 204393.6
$ws.Range("K122").Value = 6493.7145
$ws.Range("L122").Value = 613180.8
$ws.Range("M122").Value = -4043.7145
$ws.Range("N122").Value = -618080.8

$ws.Range("H123").Value = 48000
$ws.Range("J123").Value = 48000
$ws.Range("L123").Value = 48000
$ws.Range("N123").Value = -52900

$ws.Range("H132").Value = 2889.5881
$ws.Range("I132").Value = 2846.7273
$ws.Range("J132").Value = 2968.1667
$ws.Range("K132").Value = 8540.1819
$ws.Range("L132").Value = 8904.500100000001
$ws.Range("M132").Value = -6010.1819
$ws.Range("N132").Value = -13964.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3590.3333
$ws.Range("I22").Value = 2453.8696
$ws.Range("J22").Value = 10125
$ws.Range("K22").Value = 2453.8696
$ws.Range("L22").Value = 10125
$ws.Range("M22").Value = -2158.8696
$ws.Range("N22").Value = -10715

$ws.Range("H27").Value = 3590.3333
$ws.Range("I27").Value = 2453.8696
$ws.Range("J27").Value = 10125
$ws.Range("K27").Value = 2453.8696
$ws.Range("L27").Value = 10125
$ws.Range("M27").Value = -2346.8696
$ws.Range("N27").Value = -10339

$ws.Range("H42").Value = 5018010.5
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H49").Value = 5018010.5
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H55").Value = 827.4167
$ws.Range("J55").Value = 1248.75
$ws.Range("L55").Value = 1248.75
$ws.Range("N55").Value = -1594.75

$ws.Range("H68").Value = 3467.6667
$ws.Range("I68").Value = 2700
$ws.Range("K68").Value = 2700
$ws.Range("M68").Value = -1951

$ws.Range("H71").Value = 3467.6667
$ws.Range("I71").Value = 2700
$ws.Range("K71").Value = 13500
$ws.Range("M71").Value = -9756

$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988

$ws.Range("H106").Value = 20273.6
$ws.Range("J106").Value = 20273.6
$ws.Range("L106").Value = 20273.6
$ws.Range("N106").Value = -22797.6

$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060
